$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.083779228826828
$ws.Range("D2").Value = 1.031183744646027
$ws.Range("E2").Value = 1.084284111261289
$ws.Range("F2").Value = 1.090572204803973
$ws.Range("I2").Value = 1.033406222798955
$ws.Range("J2").Value = 1.088642252437694
$ws.Range("K2").Value = 1.033992598156257
$ws.Range("L2").Value = 1.086947531083654
$ws.Range("M2").Value = 1.09321941508782
$ws.Range("N2").Value = 1.090188249305167
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.087710804139586
$ws.Range("D3").Value = 1.031686969028093
$ws.Range("E3").Value = 1.087764978864458
$ws.Range("F3").Value = 1.094039887044968
$ws.Range("I3").Value = 1.03360284524045
$ws.Range("J3").Value = 1.092222853643149
$ws.Range("K3").Value = 1.034305152405643
$ws.Range("L3").Value = 1.090240982380089
$ws.Range("M3").Value = 1.096500955495478
$ws.Range("N3").Value = 1.093773935375034
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.090230394999502
$ws.Range("D4").Value = 1.032010757823463
$ws.Range("E4").Value = 1.089994480303252
$ws.Range("F4").Value = 1.096260605087055
$ws.Range("I4").Value = 1.033726680896811
$ws.Range("J4").Value = 1.094515961209113
$ws.Range("K4").Value = 1.034504852944241
$ws.Range("L4").Value = 1.092349153290535
$ws.Range("M4").Value = 1.098601157028766
$ws.Range("N4").Value = 1.096070299416766
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.091283983298007
$ws.Range("D5").Value = 1.032146448538706
$ws.Range("E5").Value = 1.090926463767977
$ws.Range("F5").Value = 1.097188833929538
$ws.Range("I5").Value = 1.03377793959412
$ws.Range("J5").Value = 1.09547446978439
$ws.Range("K5").Value = 1.034588206906088
$ws.Range("L5").Value = 1.093230109562205
$ws.Range("M5").Value = 1.099478698444767
$ws.Range("N5").Value = 1.097030169184162
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.091460559762083
$ws.Range("D6").Value = 1.032169206620338
$ws.Range("E6").Value = 1.091082642069554
$ws.Range("F6").Value = 1.097344378120144
$ws.Range("I6").Value = 1.03378649942931
$ws.Range("J6").Value = 1.095635089374806
$ws.Range("K6").Value = 1.034602167456142
$ws.Range("L6").Value = 1.09337771885123
$ws.Range("M6").Value = 1.099625730530514
$ws.Range("N6").Value = 1.09719101687282
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.09024449503567
$ws.Range("D7").Value = 1.032012572608485
$ws.Range("E7").Value = 1.090006954107997
$ws.Range("F7").Value = 1.096273028968684
$ws.Range("I7").Value = 1.033727368955074
$ws.Range("J7").Value = 1.094528790278586
$ws.Range("K7").Value = 1.034505969070712
$ws.Range("L7").Value = 1.092360945352662
$ws.Range("M7").Value = 1.098612903717239
$ws.Range("N7").Value = 1.096083146704988
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.08511309421305
$ws.Range("D8").Value = 1.031354195882566
$ws.Range("E8").Value = 1.085465323774429
$ws.Range("F8").Value = 1.091749012222105
$ws.Range("I8").Value = 1.033473380957622
$ws.Range("J8").Value = 1.089857367045092
$ws.Range("K8").Value = 1.034098759168672
$ws.Range("L8").Value = 1.088065414333026
$ws.Range("M8").Value = 1.094333326224605
$ws.Range("N8").Value = 1.091405089514683
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.075875035670154
$ws.Range("D9").Value = 1.030179641396726
$ws.Range("E9").Value = 1.077279395261191
$ws.Range("F9").Value = 1.083592280025598
$ws.Range("I9").Value = 1.0329993624412
$ws.Range("J9").Value = 1.081435279059832
$ws.Range("K9").Value = 1.033361323570496
$ws.Range("L9").Value = 1.080312944193796
$ws.Range("M9").Value = 1.086607040778298
$ws.Range("N9").Value = 1.082971041198457
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.069571646242567
$ws.Range("D10").Value = 1.029386349286363
$ws.Range("E10").Value = 1.071687556918033
$ws.Range("F10").Value = 1.078018791806455
$ws.Range("I10").Value = 1.032664864780157
$ws.Range("J10").Value = 1.075680502828729
$ws.Range("K10").Value = 1.03285574527987
$ws.Range("L10").Value = 1.075010313252785
$ws.Range("M10").Value = 1.081320687044086
$ws.Range("N10").Value = 1.077208092525024
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.066804870940165
$ws.Range("D11").Value = 1.029040278383823
$ws.Range("E11").Value = 1.069231636567411
$ws.Range("F11").Value = 1.075570581765406
$ws.Range("I11").Value = 1.032515476250015
$ws.Range("J11").Value = 1.073152610009002
$ws.Range("K11").Value = 1.032633376005647
$ws.Range("L11").Value = 1.072679764288836
$ws.Range("M11").Value = 1.078996933649694
$ws.Range("N11").Value = 1.07467660980754
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.065771284550714
$ws.Range("D12").Value = 1.028911333978897
$ws.Range("E12").Value = 1.068313959002006
$ws.Range("F12").Value = 1.074655735334717
$ws.Range("I12").Value = 1.032459288678704
$ws.Range("J12").Value = 1.072207975218208
$ws.Range("K12").Value = 1.032550247289123
$ws.Range("L12").Value = 1.071808684265176
$ws.Range("M12").Value = 1.078128342255556
$ws.Range("N12").Value = 1.073730633526992
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.065993262998279
$ws.Range("D13").Value = 1.028939011215029
$ws.Range("E13").Value = 1.068511054041554
$ws.Range("F13").Value = 1.074852224529083
$ws.Range("I13").Value = 1.032471372932197
$ws.Range("J13").Value = 1.072410862954578
$ws.Range("K13").Value = 1.032568102910486
$ws.Range("L13").Value = 1.071995782534984
$ws.Range("M13").Value = 1.078314908274802
$ws.Range("N13").Value = 1.073933809387219
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.06671955597058
$ws.Range("D14").Value = 1.029029627995708
$ws.Range("E14").Value = 1.069155893430033
$ws.Range("F14").Value = 1.075495073304819
$ws.Range("I14").Value = 1.032510846093977
$ws.Range("J14").Value = 1.073074643197845
$ws.Range("K14").Value = 1.032626515463631
$ws.Range("L14").Value = 1.07260787228169
$ws.Range("M14").Value = 1.078925248064359
$ws.Range("N14").Value = 1.074598532274568
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.067166261268918
$ws.Range("D15").Value = 1.029085406838724
$ws.Range("E15").Value = 1.069552472183638
$ws.Range("F15").Value = 1.075890421305811
$ws.Range("I15").Value = 1.032535073869484
$ws.Range("J15").Value = 1.073482862115336
$ws.Range("K15").Value = 1.032662434637913
$ws.Range("L15").Value = 1.072984276978654
$ws.Range("M15").Value = 1.079300569918088
$ws.Range("N15").Value = 1.075007330909745
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.069754458041884
$ws.Range("D16").Value = 1.02940926170044
$ws.Range("E16").Value = 1.071849798782244
$ws.Range("F16").Value = 1.078180517218019
$ws.Range("I16").Value = 1.032674682133957
$ws.Range("J16").Value = 1.075847490504377
$ws.Range("K16").Value = 1.032870429446772
$ws.Range("L16").Value = 1.075164238277954
$ws.Range("M16").Value = 1.081474155950635
$ws.Range("N16").Value = 1.077375317342325
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.071367775154603
$ws.Range("D17").Value = 1.029611710815771
$ws.Range("E17").Value = 1.073281417730293
$ws.Range("F17").Value = 1.079607539285307
$ws.Range("I17").Value = 1.032761026746357
$ws.Range("J17").Value = 1.077320939034246
$ws.Range("K17").Value = 1.032999966670659
$ws.Range("L17").Value = 1.076522280536085
$ws.Range("M17").Value = 1.082828133212495
$ws.Range("N17").Value = 1.078850858338113
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.072305209190787
$ws.Range("D18").Value = 1.02972954882415
$ws.Range("E18").Value = 1.074113133377886
$ws.Range("F18").Value = 1.080436550757971
$ws.Range("I18").Value = 1.032810952221351
$ws.Range("J18").Value = 1.078176916941275
$ws.Range("K18").Value = 1.033075191486793
$ws.Range("L18").Value = 1.077311093383019
$ws.Range("M18").Value = 1.083614550479474
$ws.Range("N18").Value = 1.079708051831953
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.072624249614349
$ws.Range("D19").Value = 1.029769686979487
$ws.Range("E19").Value = 1.074396170280739
$ws.Range("F19").Value = 1.080718661648542
$ws.Range("I19").Value = 1.032827901679518
$ws.Range("J19").Value = 1.078468203843701
$ws.Range("K19").Value = 1.03310078525162
$ws.Range("L19").Value = 1.077579503559428
$ws.Range("M19").Value = 1.083882139592121
$ws.Range("N19").Value = 1.079999752395191
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.071195054378021
$ws.Range("D20").Value = 1.029590015589864
$ws.Range("E20").Value = 1.073128164019663
$ws.Range("F20").Value = 1.079454781126263
$ws.Range("I20").Value = 1.03275180818182
$ws.Range("J20").Value = 1.077163211568434
$ws.Range("K20").Value = 1.03298610298678
$ws.Range("L20").Value = 1.076376919484812
$ws.Range("M20").Value = 1.082683210817787
$ws.Range("N20").Value = 1.078692906881209
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.066505845406333
$ws.Range("D21").Value = 1.029002954706749
$ws.Range("E21").Value = 1.068966156422868
$ws.Range("F21").Value = 1.075305923368348
$ws.Range("I21").Value = 1.032499241621056
$ws.Range("J21").Value = 1.07287933483356
$ws.Range("K21").Value = 1.032609329180193
$ws.Range("L21").Value = 1.07242777836396
$ws.Range("M21").Value = 1.078745670426884
$ws.Range("N21").Value = 1.074402946550003
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.063523375640652
$ws.Range("D22").Value = 1.028631535560414
$ws.Range("E22").Value = 1.066317740856485
$ws.Range("F22").Value = 1.072665587103774
$ws.Range("I22").Value = 1.032336396760489
$ws.Range("J22").Value = 1.07015299712144
$ws.Range("K22").Value = 1.03236935850096
$ws.Range("L22").Value = 1.069913372324495
$ws.Range("M22").Value = 1.076238353825647
$ws.Range("N22").Value = 1.071672737125617
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.065107771803543
$ws.Range("D23").Value = 1.028828655107005
$ws.Range("E23").Value = 1.067724793310947
$ws.Range("F23").Value = 1.074068373555036
$ws.Range("I23").Value = 1.032423112488318
$ws.Range("J23").Value = 1.071601484162918
$ws.Range("K23").Value = 1.032496867575093
$ws.Range("L23").Value = 1.071249364880628
$ws.Range("M23").Value = 1.077570606682377
$ws.Range("N23").Value = 1.073123281184838
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.071273110527388
$ws.Range("D24").Value = 1.029599819491675
$ws.Range("E24").Value = 1.073197423028972
$ws.Range("F24").Value = 1.079523816288073
$ws.Range("I24").Value = 1.032755975006491
$ws.Range("J24").Value = 1.077234492495398
$ws.Range("K24").Value = 1.032992368414144
$ws.Range("L24").Value = 1.07644261210066
$ws.Range("M24").Value = 1.082748705299618
$ws.Range("N24").Value = 1.078764289035265
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.078287813360789
$ws.Range("D25").Value = 1.030485056302546
$ws.Range("E25").Value = 1.079418497935246
$ws.Range("F25").Value = 1.085724040494184
$ws.Range("I25").Value = 1.033125112947471
$ws.Range("J25").Value = 1.083636366507373
$ws.Range("K25").Value = 1.033554383163982
$ws.Range("L25").Value = 1.082339973969298
$ws.Range("M25").Value = 1.088627517868472
$ws.Range("N25").Value = 1.085175254442642
